# LIVEHTA-1904: update testdata on the "prodfix" sheet.
# The scenario name/title values lose their " - 9/19/2022" date suffix,
# which also makes the old shared string unused (removed on save) and
# causes the column width next to it to shrink/bestfit to its new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prodfix")

$newName = "PRODFix_QOL_ECON - UtilityOutcome"

$ws.Range("B2").Value = $newName
$ws.Range("B5").Value = $newName
$ws.Range("B8").Value = $newName
$ws.Range("B11").Value = $newName

# Column B now holds text the same length/style as column C, so its
# width collapses to match column C's (already best-fit) width.
$ws.Range("B1").EntireColumn.ColumnWidth = 31.944010416666668

# Re-select B11 (and, by activating/selecting, drop the stale
# topLeftCell="H1" scroll position left over from the previous edit).
$ws.Activate()
$ws.Range("B11").Select()
